$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data table of updated values scraped from the commit diff.
# Each row entry lists the new Price (D), Volume 1h (E) and Hora (G) values.
# Empty string means "leave this cell untouched" (it was not changed by the commit).
$updates = @(
    @{ Row=2; D="304.18"; E="5.13%"; G="13" }
    @{ Row=3; D="35.15"; E="13.44%"; G="13" }
    @{ Row=4; D="5.179"; E="4.64%"; G="13" }
    @{ Row=5; D="0.07846"; E="6.27%"; G="13" }
    @{ Row=6; D="2.293"; E="-3.47%"; G="13" }
    @{ Row=7; D="8.049"; E="4.30%"; G="13" }
    @{ Row=8; D="3.988"; E="7.19%"; G="13" }
    @{ Row=9; D="0.9284"; E="2.10%"; G="13" }
    @{ Row=10; D="0.1003"; E="10.10%"; G="13" }
    @{ Row=11; D="0.1834"; E="8.46%"; G="13" }
    @{ Row=12; D="0.08592"; E="5.57%"; G="13" }
    @{ Row=13; D="0.03393"; E="8.68%"; G="13" }
    @{ Row=14; D="0.09931"; E="-0.36%"; G="13" }
    @{ Row=15; D="0.001496"; E="-0.21%"; G="13" }
    @{ Row=16; D="0.04646"; E="2.83%"; G="13" }
    @{ Row=17; D="0.005769"; E="-0.95%"; G="13" }
    @{ Row=18; D="3.487"; E="-0.34%"; G="13" }
    @{ Row=19; D="2.127"; E="1.37%"; G="13" }
    @{ Row=20; D="0.3420"; E="2.91%"; G="13" }
    @{ Row=21; D=""; E="2.60%"; G="13" }
    @{ Row=22; D="4.549"; E="9.62%"; G="13" }
    @{ Row=23; D="0.2380"; E="13.43%"; G="13" }
    @{ Row=24; D=""; E="1.35%"; G="13" }
    @{ Row=25; D="0.004448"; E="6.50%"; G="13" }
    @{ Row=26; D="0.0001298"; E="0.00%"; G="13" }
    @{ Row=27; D="0.0003395"; E="0.12%"; G="13" }
    @{ Row=28; D=""; E=""; G="13" }
    @{ Row=29; D=""; E=""; G="13" }
    @{ Row=30; D=""; E=""; G="13" }
    @{ Row=31; D=""; E=""; G="13" }
    @{ Row=32; D=""; E=""; G="13" }
    @{ Row=33; D=""; E=""; G="13" }
    @{ Row=34; D=""; E=""; G="13" }
    @{ Row=35; D=""; E=""; G="13" }
    @{ Row=36; D=""; E=""; G="13" }
    @{ Row=37; D=""; E=""; G="13" }
    @{ Row=38; D=""; E=""; G="13" }
    @{ Row=39; D="0.01756"; E="10.86%"; G="13" }
    @{ Row=40; D="0.04743"; E="6.14%"; G="13" }
    @{ Row=41; D="0.007690"; E="4.28%"; G="13" }
    @{ Row=42; D="0.1415"; E="6.17%"; G="13" }
    @{ Row=43; D="0.007050"; E="-25.58%"; G="13" }
    @{ Row=44; D="0.002286"; E="2.23%"; G="13" }
    @{ Row=45; D="0.009992"; E="24.52%"; G="13" }
    @{ Row=46; D="0.00005998"; E="-1.78%"; G="13" }
    @{ Row=47; D=""; E="0.05%"; G="13" }
    @{ Row=48; D="5.800"; E="126.10%"; G="13" }
    @{ Row=49; D="0.002686"; E="34.44%"; G="13" }
    @{ Row=50; D="0.00002097"; E="0.05%"; G="13" }
    @{ Row=51; D="0.0001997"; E="0.05%"; G="13" }
)

foreach ($u in $updates) {
    $r = $u.Row

    if ($u.D -ne "") {
        $cell = $ws.Cells.Item($r, 4)   # column D - Price
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }

    if ($u.E -ne "") {
        $cell = $ws.Cells.Item($r, 5)   # column E - Volume(1h)
        $cell.NumberFormat = "@"
        $cell.Value = $u.E
        $cell.Style = "Normal"
    }

    if ($u.G -ne "") {
        $cell = $ws.Cells.Item($r, 7)   # column G - Hora
        $cell.NumberFormat = "@"
        $cell.Value = $u.G
        $cell.Style = "Normal"
    }
}
